# [Silverfox] Abnormal FX color scheme
# Adds a new "대표 색상" (representative color) column M next to the existing
# Abnormal-Type table, carrying the FX color associated with each status.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column M doesn't exist yet; give every row (2-19) the same bordered cell
# formatting already used throughout the table (copied from column L) before
# filling in values.
$ws.Range("L2:L19").Copy()
$ws.Range("M2:M19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header
$ws.Cells.Item(2, 13).Value = "대표 색상"

# Representative color per status (rows without an entry stay blank, matching
# the existing sparse K/L-style columns in this sheet).
$ws.Cells.Item(18, 13).Value = "녹색 계열"    # 16 Cure
$ws.Cells.Item(12, 13).Value = "보라색 계열"  # 10 Poison
$ws.Cells.Item(5, 13).Value  = "검은색 계열"  # 03 Nearsight
$ws.Cells.Item(3, 13).Value  = "연두색 계열"  # 01 Slow
$ws.Cells.Item(9, 13).Value  = "분홍색 계열"  # 07 Charm
$ws.Cells.Item(11, 13).Value = "노란색 계열"  # 09 Divine
$ws.Cells.Item(13, 13).Value = "주황색 계열"  # 11 Flame

# Column width for the new column (~11.6 characters, matching the rest of the
# compact columns in this sheet)
$ws.Columns.Item(13).ColumnWidth = 10.86

# Matches the cursor position left by the author after finishing the edit
$ws.Range("M12").Select()
